$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 ("theta_threshold_range") is removed entirely; pie_threshold_range
# shifts up from row 6 to row 5.
$ws.Rows.Item(5).Delete()

# Updated threshold values.
$ws.Range("B2").Value = 4
$ws.Range("B3").Value = 4.0999999999999996
$ws.Range("C3").Value = 10.1

# Column C narrows (best-fit shrank once the longer "168.5" value was removed).
$ws.Columns.Item(3).ColumnWidth = 4.857142857142857

# Selection moved off the edited cell.
$ws.Range("F12").Select() | Out-Null

# Page setup was touched (paper size / orientation recorded for the sheet).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
